$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-08-30 21:01:26"
}
